# Add the new "Test Data" worksheet (becomes sheetId=3, rId3) after the
# existing "API Test cases" sheet, populate it with the dropdown / checkout
# reference data, and leave it the active sheet/tab (matches the workbook
# diff: new sheet added, activeTab moves to it, tabSelected moves off
# "UI Test cases").

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Test Data"

# Column A: dropdown sort labels
$ws.Range("A1").Value = "Dropdown Data"
$ws.Range("A2").Value = "Name (Z to A)"
$ws.Range("A3").Value = "Price (high to low)"
$ws.Range("A4").Value = "Price (low to high)"

# Column B: checkout information used by the automation suite
$ws.Range("B1").Value = "Checkout Information"
$ws.Range("B2").Value = "Dinesh"
$ws.Range("B3").Value = "Nimmala"
$ws.Range("B4").Value = 505050

# Match the authored column widths (bestFit) as closely as the host allows
$ws.Columns.Item(1).ColumnWidth = 14.751
$ws.Columns.Item(2).ColumnWidth = 18.251

# Leave the selection/active cell on B1, as in the authored sheet
$ws.Range("B1").Select()
